$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk write columns A-D and F as text values (rows 2..40) ---
$colA = @('Dynamical phase transition in assemblies of chemotactic cells','Fluid pumping and flexoelectricity can drive lumen nucleation in cell spheroids','Genuine (response-field) Langevin equations for reaction-diffusion processes ','Langevin equations for reaction-diffusion processes','Frequency regulators for the nonperturbative renormalization group in nonequilibrium systems','Fluid pumping can drive lumen nucleation and thickness oscillations in a cell spheroid','Hydraulic and electric control of cell spheroids','Hydraulic and electric control of cell spheroids','Nonlinear rheology of cellular networks','Hydraulic and electric control of a cell spheroid','Stochastic dynamics of chemotactic colonies with logistic growth','Stochastic dynamics of chemotactic colonies with logistic growth','Hydraulic and electric properties of tissues','Nonequilibrium critical phenomena: exact Langevin equations, erosion of tilted landscapes','Collective dynamics of chemotactic cells','Collective dynamics of chemotactic cells','Collective dynamics of chemotactic cells','Quorum effects in assemblies of chemotactic cells','Polarization mechanism for chemotaxis and exact scaling exponent in assemblies of cells','Fluid pumping, lumen nucleation and electro-hydraulic phenomena in cell assemblies','Fluid pumping, lumen nucleation and electro-hydraulic phenomena in cell assemblies','Active mechanical and electrohydraulic properties of tissues','Active mechanical and electrohydraulic properties of tissues','Active mechanical and electrohydraulic properties of tissues','Active mechanical and electrohydraulic properties of tissues','Renormalization group approach to the collective dynamics of chemotactic cells','Hysteresis, phase transitions and dangerous transients in power distribution systems','Langevin equations for reaction-diffusion processes','Langevin equations for reaction-diffusion processes','Fluid pumping and active flexoelectricity can promote lumen nucleation in cell assemblies','Fluid pumping and active flexoelectricity can promote lumen nucleation in cell assemblies','Fluid pumping and active flexoelectricity can promote lumen nucleation in cell assemblies','Fluid pumping and active flexoelectricity can promote lumen nucleation in cell assemblies','Fluid pumping and active flexoelectricity can promote lumen nucleation in cell assemblies','Hydraulic and electric control of cell spheroids','Hydraulic and electric control of cell spheroids','Electrohydraulics of cells and tissues','Electrohydraulics of cells and tissues','Stochastic dynamics of chemotactic colonies with logistic growth')
$colB = @('Invited talk','Invited talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Contributed talk','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Seminar','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster','Poster')
$colC = @('inv-talk-1','inv-talk-2','talk-1','talk-2','talk-3','talk-4','talk-5','talk-6','talk-7','talk-8','talk-9','talk-10','talk-11','seminar-1','seminar-2','seminar-3','seminar-4','seminar-5','seminar-6','seminar-7','seminar-8','seminar-9','seminar-10','seminar-11','seminar-12','seminar-13','poster-1','poster-2','poster-3','poster-4','poster-5','poster-6','poster-7','poster-8','poster-9','poster-10','poster-11','poster-12','poster-13')
$colD = @('DPG Spring Meeting, Focus Session on Stochastic processes in biology','Labex "Who am I?" 2022 annual meeting','Journées de Physique Statistique','DPG Spring Meeting','Exact Renormalization Group (ERG2018)','Circle Meeting','Virtual APS March Meeting','Virtual DPG Spring Meeting','Physical Biology Circle Meeting','New frontiers in liquid matter','Journées de Physique Statistique','Sitges conference on Statistical Mechanics','Physics meets Biology','Max Planck Institute for the Physics of Complex Systems','Laboratoire de Physique','Laboratoire Interdisciplinaire de Physique','Laboratoire Matière et Systèmes Complexes','Laboratoire de Physique Théorique et Modèles Statistiques (online seminar)','Centre de Physique Théorique (online seminar)','Laboratoire Interdisciplinaire de Physique (online seminar)','Laboratoire Jean Perrin (online seminar)','Laboratoire Matière et Systèmes Complexes','Laboratoire Physico-Chimie Curie','Laboratoire Physique et Mécanique des Milieux Hétérogènes','Laboratoire Jean Perrin ','Laboratoire de Physique Théorique de la Matière Condensée','Optimization and Control of Smart Grids','Beg Rohu summer school','42nd Conference of the Middle European Cooperation in Statistical Physics','Les Houches summer school','Key Challenges in Statistical Physics (MECO 44)','Mechanical Forces in Development','Fluid Physics of Life','Organoids : Modelling Organ Development and Disease in 3D Culture','EMBO Workshop Physics of living systems: From molecules to tissues','Physics Meets Biology','Institute Scientific Evaluation','Engineering Life – Active Matter Across Scales','New Perspectives in Active Systems')
$colF = @('Regensburg Universität,  Regensburg, Germany','ENSA, Paris, France','ESPCI, Paris, France','Technische Universität Berlin, Berlin, Germany','Sorbonne Université, Paris, France','Saarland Universität, Saarbrücken, Germany','online conference','online conference','Max Planck Institute for the Physics of Complex Systems, Dresden, Germany','Sorbonne Université, Paris, France','École Normale Supérieure de Paris, France','Hotel Calipolis, Sitges, Spain','Rice Global Paris Center, Paris, France','Dresden, Germany','École Normale Supérieure de Lyon, Lyon, France','Université Grenoble-Alpes, Grenoble, France','Université de Paris, Paris, France','Université Paris-Orsay, Orsay, France','Aix-Marseille Université, Marseille, France','Université Grenoble-Alpes, Grenoble, France','Sorbonne Université, Paris, France ','Université de Paris, Paris, France','Institut Curie, Paris, France','ESPCI, Paris, France','Sorbonne Université, Paris, France ','Sorbonne Université, Paris, France ','Santa Fe, New Mexico, USA','Beg Rohu, France','École Normale Supérieure de Lyon, Lyon, France','Les Houches, France','Kloster Seeon, Munich, Germany','EMBL, Heidelberg, Germany','Max Planck Institute for the Physics of Complex Systems, Dresden, Germany','online conference','online conference','online conference','Max Planck Institute for the Physics of Complex Systems, Dresden, Germany','Steigenberger Hotel de Saxe Dresden, Germany','Max Planck Institute for the Physics of Complex Systems, Dresden, Germany')
$colE = @(43559,44907,42397,43171,43291,43551,44271,44278,44454,44749,44952,45075,45223,42824,43837,43839,43845,43921,44118,44137,44145,44543,44589,44608,44614,44844,41050,42583,42775,43313,43586,43649,43759,44125,44354,44403,44656,44740,45041)
$cstyles = @(0,0,0,0,0,2,2,2,2,2,2,2,2,0,0,0,0,2,2,2,2,2,2,2,2,2,0,0,0,0,0,0,0,0,0,0,0,0,0)
$estyles = @(1,1,1,1,1,1,1,1,3,1,1,1,1,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,1,1,1)

$n = $colA.Length
$arrA = New-Object "object[,]" $n,1
$arrB = New-Object "object[,]" $n,1
$arrC = New-Object "object[,]" $n,1
$arrD = New-Object "object[,]" $n,1
$arrE = New-Object "object[,]" $n,1
$arrF = New-Object "object[,]" $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arrA[$i,0] = $colA[$i]
    $arrB[$i,0] = $colB[$i]
    $arrC[$i,0] = $colC[$i]
    $arrD[$i,0] = $colD[$i]
    $arrE[$i,0] = $colE[$i]
    $arrF[$i,0] = $colF[$i]
}

$ws.Range("A2:A40").Value = $arrA
$ws.Range("B2:B40").Value = $arrB
$ws.Range("C2:C40").Value = $arrC
$ws.Range("D2:D40").Value = $arrD
$ws.Range("E2:E40").Value = $arrE
$ws.Range("F2:F40").Value = $arrF
$ws.Range("G2:G40").ClearContents()

Write-Output "values-written"
